# Applies the "论文总结.xlsx" update:
#  - F2 gets replaced with the new 3-part summary text (was the now-removed
#    "基于规则与语义编辑距离的简称" string) and gets wrap-text turned on.
#  - Four new rows (13-16) of paper notes are appended.
#  - Column A is widened and the sheet selection/frozen-pane view is moved
#    down to show the newly added rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- F2: replace the short placeholder with the fuller breakdown ------------
$ws.Range("F2").Value = "主要分成三个部分：`n1.基于条件随机场的全称识别算法`n2.基于规则与语义编辑距离的简称识别`n3.基于全称与简称的命名实体消歧算法"
$ws.Range("F2").WrapText = $true

# --- Row 13: SparkCRF paper ---------------------------------------------
$ws.Range("A13").Value = "SparkCRF：一种基于Spark的并行CRFs算法实现"
$ws.Range("B13").Value = 2016.3
$ws.Range("C13").Value = "期刊"
$ws.Range("D13").Value = "设计了SparkCRF"

# --- Row 14: Random Walks NED paper --------------------------------------
$ws.Range("A14").Value = "Robust Named Entity Disambiguation withRandom Walks"
$ws.Range("E14").Value = "在自然语言处理中，关于命名体有两个重要的任务，一个是命名体的识别(Named Entity Recognition)，另一个是命名体的消歧(Named Entity Disambiguation)"
$ws.Range("E14").WrapText = $true
$ws.Rows.Item(14).RowHeight = 42.75

# --- Row 15: Bidirectional LSTM-CRF paper --------------------------------
$ws.Range("A15").Value = "Bidirectional LSTM-CRF models for sequence`ntagging"
$ws.Range("A15").WrapText = $true
$ws.Range("E15").Value = "与传统的神经网络相比，RNN多了上一层隐藏层与当前层隐藏层的链接，因此可以用来记忆历史的信息`nLSTM与RNN类似，只是隐藏层的结构发生了变化，因此LSTM可以更好地寻找与挖掘长距离的依赖。`n"
$ws.Range("E15").WrapText = $true
$ws.Rows.Item(15).RowHeight = 71.25

# --- Row 16: financial-domain event extraction note ----------------------
$ws.Range("A16").Value = "金融领域的事件句抽取"
$ws.Range("E16").Value = "关于简称的抽取与识别可以参考一下"
$ws.Range("E16").WrapText = $true

# --- Column / view tweaks -------------------------------------------------
# Target stored width is 41.625 "characters"; the engine's column-width grid
# is quantized to 1/7ths of a character (Excel's pixel-rounding behaviour),
# so 40.85 is the closest achievable ColumnWidth input (lands on 41.571...).
$ws.Columns.Item(1).ColumnWidth = 40.85

$ws.Range("E16").Select() | Out-Null
